# Rachel_money.xlsx edit script
# Commit message: "Adds field in career stats to bound prediction"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Remove the unused empty "Sheet1" worksheet
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2. Career Stats: update a couple of values and add a new "Project Through"
#    field used to bound the projection.
# ---------------------------------------------------------------------------
$stats = $wb.Worksheets.Item("Career Stats")
$stats.Range("B7").Value = 0
$stats.Range("B8").Value = 33000
$stats.Range("A11").Value = "Project Through"
$stats.Range("B11").Value = 2033

# ---------------------------------------------------------------------------
# 3. Career Projection: revise existing projection rows and add new ones
# ---------------------------------------------------------------------------
$proj = $wb.Worksheets.Item("Career Projection")

# Row 2 updates
$proj.Range("C2").Value = 43852
$proj.Range("E2").ClearContents() | Out-Null
$proj.Range("F2").ClearContents() | Out-Null
$proj.Range("G2").Value = 45060
$proj.Range("H2").Value = 43000
$proj.Range("H2").NumberFormat = "`"$`"#,##0_);[Red]\(`"$`"#,##0\)"

# Row 3 updates
$proj.Range("C3").Value = 45060
$proj.Range("D3").Value = 20755
$proj.Range("G3").Value = 43831
$proj.Range("J3").Value = 2000

# Row 4 previously only had G4/I4 - clear it out, it is replaced below
$proj.Range("G4").ClearContents() | Out-Null
$proj.Range("I4").ClearContents() | Out-Null

# New row 4 - copy date/number formats from row 2's analogous cells first
$proj.Range("A2").Copy() | Out-Null
$proj.Range("A4").PasteSpecial(-4122) | Out-Null
$proj.Range("C2").Copy() | Out-Null
$proj.Range("C4").PasteSpecial(-4122) | Out-Null
$proj.Range("D2").Copy() | Out-Null
$proj.Range("D4").PasteSpecial(-4122) | Out-Null
$proj.Range("G2").Copy() | Out-Null
$proj.Range("G4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$proj.Range("A4").Value = 47993
$proj.Range("B4").Value = "O-5"
$proj.Range("C4").Value = 46138
$proj.Range("D4").Value = 80841
$proj.Range("G4").Value = 46138
$proj.Range("H4").Value = 38000
$proj.Range("H4").NumberFormat = "`"$`"#,##0_);[Red]\(`"$`"#,##0\)"

# New row 5
$proj.Range("C2").Copy() | Out-Null
$proj.Range("C5").PasteSpecial(-4122) | Out-Null
$proj.Range("D2").Copy() | Out-Null
$proj.Range("D5").PasteSpecial(-4122) | Out-Null
$proj.Range("G2").Copy() | Out-Null
$proj.Range("G5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$proj.Range("C5").Value = 47279
$proj.Range("D5").Value = 78251
$proj.Range("G5").Value = 47289
$proj.Range("H5").Value = 40000

# New row 6
$proj.Range("G2").Copy() | Out-Null
$proj.Range("G6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$proj.Range("G6").Value = 44197
$proj.Range("J6").Value = 0

# ---------------------------------------------------------------------------
# 4. Assets: re-prioritize / rename rows, add a new "Checking" row
# ---------------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")
$assets.Range("C3").Value = 4000
$assets.Range("B4").Value = "Brokerage"
$assets.Range("C4").Value = 0
$assets.Range("D4").Value = 6
$assets.Range("B5").Value = "Savings"
$assets.Range("D5").Value = 1.45

$assets.Range("A5:D5").Copy() | Out-Null
$assets.Range("A6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$assets.Range("A6").Value = 5
$assets.Range("B6").Value = "Checking"
$assets.Range("C6").Value = 4000
$assets.Range("D6").Value = 0

# ---------------------------------------------------------------------------
# 5. Debts: add a leading "Priority" column ahead of the existing headers
# ---------------------------------------------------------------------------
$debts = $wb.Worksheets.Item("Debts")
$debtName = $debts.Range("A1").Value()
$debtBalance = $debts.Range("B1").Value()
$interestRate = $debts.Range("C1").Value()

# Reuse the bold header formatting already used on the Assets sheet
$assets.Range("A1:D1").Copy() | Out-Null
$debts.Range("A1:D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$debts.Range("D1").Value = $interestRate
$debts.Range("C1").Value = $debtBalance
$debts.Range("B1").Value = $debtName
$debts.Range("A1").Value = "Priority"

# ---------------------------------------------------------------------------
# 6. Make "Career Stats" the active sheet (mirrors the authored workbook view)
# ---------------------------------------------------------------------------
$stats.Activate() | Out-Null
